$d = $word.ActiveDocument

# Before the edit, the document body holds three paragraphs:
#   1) the bookmark paragraph (bookmark "_GoBack") that must be kept as-is
#   2) an empty paragraph that only carries a tab stop / first-line indent
#   3) a further empty trailing paragraph that owns the section properties
#      (<w:sectPr>) for the document
#
# The edit removes paragraphs 2 and 3 entirely, so that the section
# properties end up attached directly after the bookmark paragraph with no
# intervening empty paragraphs - matching how the document looked before
# those two stray empty paragraphs were introduced.

$paras = $d.Paragraphs

# Step 1: delete the 2nd paragraph (the one with the tab stop / indent).
# Deleting its paragraph mark merges it away; the former 3rd (final,
# sectPr-owning) paragraph becomes the new 2nd paragraph.
$paras.Item(2).Range.Delete()

# Step 2: delete the whole document story. Because the range spans both the
# (kept) bookmark paragraph's own mark and the still-remaining trailing
# empty paragraph, this collapses the document down to a single paragraph -
# the original bookmark paragraph - immediately followed by the section
# properties, with its content (the bookmark) preserved.
$d.Range(0, $d.Content.End).Delete()
